# This edit inserts a brand-new data row at row 422 (pushing the previously
# existing rows 422-508 down to 423-509), and fills the new row 422 with the
# new record's values (a "$/paquete" priced Albahaca entry for the Región de
# Arica y Parinacota on 2022-08-09 / serial date 44782).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 422 and below down by one row.
$ws.Rows.Item(422).Insert()

# Populate the newly inserted row 422 with its values.
$ws.Cells.Item(422, 1).Value  = 6
$ws.Cells.Item(422, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(422, 3).Value  = "Metropolitana"
$ws.Cells.Item(422, 4).Value  = 44782
$ws.Cells.Item(422, 5).Value  = 13
$ws.Cells.Item(422, 6).Value  = 100112052
$ws.Cells.Item(422, 7).Value  = "Albahaca"
$ws.Cells.Item(422, 8).Value  = "Sin especificar"
$ws.Cells.Item(422, 9).Value  = "Primera"
$ws.Cells.Item(422, 10).Value = 300
$ws.Cells.Item(422, 11).Value = 4000
$ws.Cells.Item(422, 12).Value = 4500
$ws.Cells.Item(422, 13).Value = 4250
$ws.Cells.Item(422, 14).Value = "`$/paquete"
$ws.Cells.Item(422, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(422, 16).Value = 4250
$ws.Cells.Item(422, 17).Value = 1
$ws.Cells.Item(422, 18).Value = "Hortaliza"
